# edit.ps1
# Bug fix: senior editor schedule generator no longer asks for display name.
# This updates the header names for the two senior editor columns (H, I)
# and fixes the shift values that were generated incorrectly for the
# senior editors (now uniformly "13-22" where applicable) as well as a
# handful of regular editor shift corrections that resulted from
# regenerating the schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("H1").Value = "Daisy"
    $ws.Range("I1").Value = "Tak"
    $ws.Range("B2").Value = "15-24"
    $ws.Range("D2").Value = "7-16"
    $ws.Range("H2").Value = "13-22"
    $ws.Range("I2").Value = "13-22"
    $ws.Range("H3").Value = "13-22"
    $ws.Range("I3").Value = "13-22"
    $ws.Range("B4").Value = "7-16"
    $ws.Range("F4").Value = "off"
    $ws.Range("G4").Value = "10-19"
    $ws.Range("H4").Value = "13-22"
    $ws.Range("I4").Value = "13-22"
    $ws.Range("C5").Value = "7-16"
    $ws.Range("D5").Value = "10-19"
    $ws.Range("H5").Value = "13-22"
    $ws.Range("I5").Value = "13-22"
    $ws.Range("C6").Value = "7-16"
    $ws.Range("H6").Value = "13-22"
    $ws.Range("I6").Value = "13-22"
    $ws.Range("B7").Value = "15-24"
    $ws.Range("H7").Value = "13-22"
    $ws.Range("I7").Value = "13-22"
    $ws.Range("D8").Value = "15-24"
    $ws.Range("E8").Value = "7-16"
    $ws.Range("F8").Value = "10-19"
    $ws.Range("G8").Value = "off"
    $ws.Range("H8").Value = "13-22"
    $ws.Range("I8").Value = "13-22"
    $ws.Range("C9").Value = "10-19"
    $ws.Range("D9").Value = "7-16"
    $ws.Range("H9").Value = "13-22"
    $ws.Range("I9").Value = "13-22"
    $ws.Range("G10").Value = "7-16"
    $ws.Range("H10").Value = "13-22"
    $ws.Range("I10").Value = "13-22"
    $ws.Range("G11").Value = "10-19"
    $ws.Range("H11").Value = "13-22"
    $ws.Range("I11").Value = "13-22"
    $ws.Range("E12").Value = "15-24"
    $ws.Range("F12").Value = "10-19"
    $ws.Range("G12").Value = "7-16"
    $ws.Range("H12").Value = "13-22"
    $ws.Range("I12").Value = "13-22"
    $ws.Range("H13").Value = "13-22"
    $ws.Range("I13").Value = "13-22"
    $ws.Range("H14").Value = "13-22"
    $ws.Range("I14").Value = "13-22"
    $ws.Range("F15").Value = "7-16"
    $ws.Range("H15").Value = "13-22"
    $ws.Range("I15").Value = "13-22"
    $ws.Range("G16").Value = "7-16"
    $ws.Range("H16").Value = "13-22"
    $ws.Range("I16").Value = "13-22"
    $ws.Range("H17").Value = "7-16"
    $ws.Range("I17").Value = "7-16"
    $ws.Range("H18").Value = "13-22"
    $ws.Range("I18").Value = "13-22"
    $ws.Range("D19").Value = "10-19"
    $ws.Range("E19").Value = "off"
    $ws.Range("F19").Value = "15-24"
    $ws.Range("H19").Value = "13-22"
    $ws.Range("I19").Value = "13-22"
    $ws.Range("C20").Value = "10-19"
    $ws.Range("D20").Value = "off"
    $ws.Range("E20").Value = "15-24"
    $ws.Range("F20").Value = "off"
    $ws.Range("H20").Value = "13-22"
    $ws.Range("I20").Value = "13-22"
    $ws.Range("B21").Value = "7-16"
    $ws.Range("E21").Value = "10-19"
    $ws.Range("F21").Value = "off"
    $ws.Range("H21").Value = "13-22"
    $ws.Range("I21").Value = "13-22"
    $ws.Range("H22").Value = "13-22"
    $ws.Range("I22").Value = "13-22"
    $ws.Range("H23").Value = "13-22"
    $ws.Range("I23").Value = "13-22"
    $ws.Range("H24").Value = "13-22"
    $ws.Range("I24").Value = "13-22"
    $ws.Range("H25").Value = "13-22"
    $ws.Range("I25").Value = "13-22"
    $ws.Range("H26").Value = "13-22"
    $ws.Range("I26").Value = "13-22"
    $ws.Range("E27").Value = "7-16"
    $ws.Range("G27").Value = "off"
    $ws.Range("H27").Value = "13-22"
    $ws.Range("I27").Value = "13-22"
    $ws.Range("B28").Value = "10-19"
    $ws.Range("G28").Value = "off"
    $ws.Range("H28").Value = "13-22"
    $ws.Range("I28").Value = "13-22"
    $ws.Range("H29").Value = "13-22"
    $ws.Range("I29").Value = "13-22"
    $ws.Range("H30").Value = "13-22"
    $ws.Range("I30").Value = "13-22"
    $ws.Range("H31").Value = "13-22"
    $ws.Range("I31").Value = "13-22"
